$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Metadata sheet updates
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B3").Value  = "0.4.0-snapshot-1"                 # Version
$wsMeta.Range("B6").Value  = "draft"                             # Status
$wsMeta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"        # Date
$wsMeta.Range("B10").Value = "ANS (https://esante.gouv.fr)"     # Contact

# ---------------------------------------------------------------------
# 2) Elements sheet updates: swap the "Mapping: RIM Mapping" (col AK/37)
#    and "Mapping: Spécification métier ..." (col AL/38) columns,
#    including their header text, their data, and their widths.
# ---------------------------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Header row (row 1)
$wsElem.Range("AK1").Value = "Mapping: Spécification métier vers l'extension ROR AvailableTimeEffectiveOpeningClosingDate"
$wsElem.Range("AL1").Value = "Mapping: RIM Mapping"

# Data rows 2-10 (swap AK <-> AL content)
$akVals = @{}
$alVals = @{}
for ($r = 2; $r -le 10; $r++) {
    $akVals[$r] = $wsElem.Cells.Item($r, 37).Text
    $alVals[$r] = $wsElem.Cells.Item($r, 38).Text
}
for ($r = 2; $r -le 10; $r++) {
    $wsElem.Cells.Item($r, 37).Value = $alVals[$r]
    $wsElem.Cells.Item($r, 38).Value = $akVals[$r]
}

# Column widths: swap widths of columns AK (37) and AL (38)
# (AK becomes the wide column ~99.84 chars, AL becomes the narrow column ~24.98 chars)
$wsElem.Columns.Item(37).ColumnWidth = 99.0
$wsElem.Columns.Item(38).ColumnWidth = 24.166666666666664
